# feat: add 2022-Q1 data
#
# 1) Insert a new worksheet "2022-Q1" right after "2021-Q4" (before "总计")
#    holding the per-fund holding detail for the new quarter.
# 2) Insert a new first data row into "总计" summarising that quarter
#    (2022-Q1, 7 funds held, 8.53 亿元 total market value), pushing the
#    existing rows down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Build the new "2022-Q1" sheet, placed right after "2021-Q4".
# ---------------------------------------------------------------------
$after = $wb.Worksheets.Item("2021-Q4")
$q1 = $wb.Worksheets.Add($null, $after)
$q1.Name = "2022-Q1"

# Copy formatting (fonts/borders/alignment) from the sibling "2021-Q4"
# sheet so the new sheet's header row / index column match the existing
# bold-bordered-centered style used across the other per-quarter sheets.
$after.Range("A1:H8").Copy()
$q1.Range("A1").PasteSpecial(-4122)
$q1.Range("A1").ClearContents()

$headers = @("", "基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($c = 2; $c -le 8; $c++) {
    $q1.Cells.Item(1, $c).Value = $headers[$c - 1]
}

# Column B (基金代码) has leading zeros, and columns D (基金规模),
# E (股票总仓位), F (仓位占比), G (持有市值) keep trailing zeros (e.g.
# "4.7374" / "0.2910") — all of these are stored as text in the source
# data, so force a text number format before assigning values, otherwise
# Excel's autodetect would coerce them to numbers and mangle them.
$q1.Range("B2:B8").NumberFormat = "@"
$q1.Range("D2:G8").NumberFormat = "@"

$rows = @(
    @(0, "011363", "南方兴润价值一年持有期混合A", "123.05", "62.78", "3.85", "4.7374", 5),
    @(1, "202003", "南方绩优成长混合A",           "48.96",  "65.45", "4.68", "2.2913", 1),
    @(2, "011364", "南方兴润价值一年持有期混合C", "19.67",  "62.78", "3.85", "0.7573", 5),
    @(3, "012412", "汇泉策略优选混合型证券投资基金", "23.52", "70.35", "1.82", "0.4281", 10),
    @(4, "501062", "南方瑞合三年定期开放混合(LOF)", "6.88",  "56.71", "4.23", "0.2910", 6),
    @(5, "006540", "南方绩优成长混合C",           "0.38",   "65.45", "4.68", "0.0178", 1),
    @(6, "320016", "诺安多策略混合",               "0.19",   "80.02", "4.78", "0.0091", 3)
)

$r = 2
foreach ($row in $rows) {
    $q1.Cells.Item($r, 1).Value = $row[0]
    $q1.Cells.Item($r, 2).Value = $row[1]
    $q1.Cells.Item($r, 3).Value = $row[2]
    $q1.Cells.Item($r, 4).Value = $row[3]
    $q1.Cells.Item($r, 5).Value = $row[4]
    $q1.Cells.Item($r, 6).Value = $row[5]
    $q1.Cells.Item($r, 7).Value = $row[6]
    $q1.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 2. Prepend the 2022-Q1 summary row to the "总计" sheet.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()

# The inserted row picks up formatting from the row above (the header);
# clear it and reapply just the index-column style (matches A3:A6).
$total.Range("A2:D2").ClearFormats()
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q1"
$total.Cells.Item(2, 3).Value = 7
$total.Cells.Item(2, 4).Value = 8.529999999999999

# The index column (A) is a plain 0-based row counter, not a carried-over
# value, so renumber the rows that got pushed down (previously 0..3,
# now rows 3..6 and renumbered 1..4).
$total.Cells.Item(3, 1).Value = 1
$total.Cells.Item(4, 1).Value = 2
$total.Cells.Item(5, 1).Value = 3
$total.Cells.Item(6, 1).Value = 4
